# guide41_login.xlsx — apply the recorded edits via Excel COM automation.
#
# Summary of the change:
#  - "index" sheet becomes the active/selected tab (cursor on B6).
#  - "p2" sheet: cursor moves to B15; the login-screenshot rows (7-9) are
#    edited - row 7 loses its screenshot reference, row 8 now points at
#    login1.png, and a new row 9 picks up a chartn/login2.png pair; also
#    gets an explicit Page Setup (A4/portrait).
#  - "p3" sheet: cursor moves to B13 and it is no longer the active tab
#    (index takes over as the active tab instead).

$wb = $excel.ActiveWorkbook

$wsIndex = $wb.Worksheets.Item("index")
$wsP1    = $wb.Worksheets.Item("p1")
$wsP2    = $wb.Worksheets.Item("p2")
$wsP3    = $wb.Worksheets.Item("p3")

# ---------------------------------------------------------------------
# p2 ("p2" sheet): rework the screenshot table rows 7-9.
# ---------------------------------------------------------------------
# Populate the new row 9 (C9/D9) first, copying formatting from row 8,
# while the shared strings "login1.png" / "login2.png" are still
# referenced elsewhere - this keeps the shared-string table stable.
$wsP2.Range("C8").Copy()
$wsP2.Range("C9").PasteSpecial(-4122)   # xlPasteFormats
$wsP2.Range("C9").Value = "chartn"

$wsP2.Range("D8").Copy()
$wsP2.Range("D9").PasteSpecial(-4122)   # xlPasteFormats
$wsP2.Range("D9").Value = "login2.png"

# Row 8's screenshot becomes the one that used to live in row 7.
$wsP2.Range("D8").Value = "login1.png"

# Row 7 no longer references any screenshot: clear D7 completely (it
# disappears from the row) and blank out C7 (keeps its style, loses its
# value).
$wsP2.Range("D7").Clear()
$wsP2.Range("C7").ClearContents()

# p2 gets an explicit page setup (A4 portrait).
$wsP2.PageSetup.PaperSize = 9      # xlPaperA4
$wsP2.PageSetup.Orientation = 1    # xlPortrait

# ---------------------------------------------------------------------
# Selection / active-tab bookkeeping.
# Select on the non-final sheets first, finish on "index" so it ends up
# as the active tab/sheet when the workbook is saved.
# ---------------------------------------------------------------------
$wsP2.Range("B15").Select()
$wsP3.Range("B13").Select()
$wsIndex.Range("B6").Select()
$wsIndex.Activate()
